# The document's single section has a "first page" header/footer pair
# (Headers/Footers Item(2), i.e. wdHeaderFooterFirstPage) and a "default"
# header/footer pair (Headers/Footers Item(1), i.e. wdHeaderFooterPrimary).
# Each of the four headers/footers carries exactly one inline picture:
#   - the headers show the BTEC logo (AlternativeText "BTec_Logo-Orange")
#   - the footers show the Pearson logo (AlternativeText ends in
#     "PearsonLogo.png")
# The pictures' internal Word "name" metadata got swapped/mislabeled the
# last time they were pasted in, so fix it up:
#   headers (BTEC logo):    image2.jpg -> image1.jpg
#   footers (Pearson logo): image1.png -> image2.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

Write-Host "Fixing inline picture names in headers/footers..."

for ($hi = 1; $hi -le 2; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shape = $shapes.Item($i)
            Write-Host "Header $hi shape $i AlternativeText:" $shape.AlternativeText
            if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                $shape.Name = "image1.jpg"
                Write-Host "Header $hi shape $i renamed to image1.jpg"
            }
        }
    }
}

for ($fi = 1; $fi -le 2; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shape = $shapes.Item($i)
            Write-Host "Footer $fi shape $i AlternativeText:" $shape.AlternativeText
            if ($shape.AlternativeText -like "*PearsonLogo.png") {
                $shape.Name = "image2.png"
                Write-Host "Footer $fi shape $i renamed to image2.png"
            }
        }
    }
}

Write-Host "Done."
